# Generate Report for Handoff
# A new handoff was generated for b.md (zh-cn and de-de), so update the
# localization-status report to reflect the new "Ready for handoff" state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: update the summary row for b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-12 18:44:47"

# ---------------------------------------------------------------------
# zh-cn sheet: update the detail row for b.md (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-12 18:44:39"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/a4d40b43bf6ae511fce17018876cbc5249960135/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e9207a6b7b86df2f36e8be9a7dfc5f0c22825e4a/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet: update the detail row for b.md (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-12 18:44:47"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/a4d40b43bf6ae511fce17018876cbc5249960135/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e9207a6b7b86df2f36e8be9a7dfc5f0c22825e4a/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 40
